$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 136.9958801269531
$ws.Range("B3").Value = 138.0481414794922
$ws.Range("B4").Value = 140.4620361328125
$ws.Range("B5").Value = 136.2327423095703
$ws.Range("B6").Value = 137.4442443847656
$ws.Range("B7").Value = 137.7482147216797
$ws.Range("B8").Value = 140.3592529296875
$ws.Range("B9").Value = 138.2991180419922
$ws.Range("B10").Value = 138.9996643066406
$ws.Range("B11").Value = 136.6790771484375
$ws.Range("B12").Value = 137.0670013427734
$ws.Range("B13").Value = 141.1650543212891
$ws.Range("B14").Value = 142.7672882080078
$ws.Range("B15").Value = 151.2744903564453
$ws.Range("B16").Value = 154.7579650878906
$ws.Range("B17").Value = 190.2874145507812
$ws.Range("B18").Value = 178.5078277587891
$ws.Range("B19").Value = 185.5553436279297
$ws.Range("B20").Value = 175.2987823486328
$ws.Range("B21").Value = 177.3972015380859
$ws.Range("B22").Value = 178.0746307373047
$ws.Range("B23").Value = 174.8529968261719
$ws.Range("B24").Value = 173.6699981689453
$ws.Range("B25").Value = 175.040283203125
$ws.Range("B26").Value = 176.373046875
$ws.Range("B27").Value = 175.2969360351562
$ws.Range("B28").Value = 180.6457824707031
$ws.Range("B29").Value = 169.285888671875
$ws.Range("B30").Value = 173.5343933105469
$ws.Range("B31").Value = 175.0888061523438
$ws.Range("B32").Value = 181.9594421386719
$ws.Range("B33").Value = 201.9370422363281
$ws.Range("B34").Value = 193.9628601074219
$ws.Range("B35").Value = 233.6648254394531
$ws.Range("B36").Value = 244.1379089355469
$ws.Range("B37").Value = 243.9602966308594
$ws.Range("B38").Value = 229.3142547607422
$ws.Range("B39").Value = 216.3759918212891
$ws.Range("B40").Value = 189.3650970458984
$ws.Range("B41").Value = 178.8901977539062
$ws.Range("B42").Value = 171.4137573242188
$ws.Range("B43").Value = 151.23828125
$ws.Range("B44").Value = 155.5109100341797
$ws.Range("B45").Value = 134.5513763427734
$ws.Range("B46").Value = 141.1094818115234
$ws.Range("B47").Value = 124.3529891967773
$ws.Range("B48").Value = 134.5599975585938
$ws.Range("B49").Value = 121.3530883789062
